# Update countries & provincias Spain
# Refresh COVID stats and re-sort a couple of neighboring rows whose
# "Casos totales" values now land in a different relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Rusia - updated stats
$ws.Cells.Item(7, 2).Value = 1204502
$ws.Cells.Item(7, 3).Value = 9859
$ws.Cells.Item(7, 4).Value = 975859
$ws.Cells.Item(7, 5).Value = 207392
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 174
$ws.Cells.Item(7, 8).Value = 21251

# Row 23: Filipinas - updated stats
$ws.Cells.Item(23, 2).Value = 319330
$ws.Cells.Item(23, 3).Value = 2674
$ws.Cells.Item(23, 4).Value = 255046
$ws.Cells.Item(23, 5).Value = 58606
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 62
$ws.Cells.Item(23, 8).Value = 5678

# Rows 43-44: Polonia moves above Emiratos Arabes Unidos (re-sorted by Casos totales)
$ws.Cells.Item(43, 1).Value = "Polonia"
$ws.Cells.Item(43, 2).Value = 98140
$ws.Cells.Item(43, 3).Value = 2367
$ws.Cells.Item(43, 4).Value = 72209
$ws.Cells.Item(43, 5).Value = 23327
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 34
$ws.Cells.Item(43, 8).Value = 2604

$ws.Cells.Item(44, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(44, 2).Value = 96529
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = 86071
$ws.Cells.Item(44, 5).Value = 10034
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 424

# Row 59: Singapur - updated stats
$ws.Cells.Item(59, 2).Value = 57800
$ws.Cells.Item(59, 3).Value = 6
$ws.Cells.Item(59, 4).Value = 57534
$ws.Cells.Item(59, 5).Value = 239
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 27

# Row 89: Croacia - updated stats
$ws.Cells.Item(89, 2).Value = 17401
$ws.Cells.Item(89, 3).Value = 241
$ws.Cells.Item(89, 4).Value = 15661
$ws.Cells.Item(89, 5).Value = 1447
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 293

# Rows 98-100: Eslovaquia moves above Malasia and Montenegro (re-sorted by Casos totales)
$ws.Cells.Item(98, 1).Value = "Eslovaquia"
$ws.Cells.Item(98, 2).Value = 12321
$ws.Cells.Item(98, 3).Value = 704
$ws.Cells.Item(98, 4).Value = 4793
$ws.Cells.Item(98, 5).Value = 7474
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 54

$ws.Cells.Item(99, 1).Value = "Malasia"
$ws.Cells.Item(99, 2).Value = 11771
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 10095
$ws.Cells.Item(99, 5).Value = 1540
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 136

$ws.Cells.Item(100, 1).Value = "Montenegro"
$ws.Cells.Item(100, 2).Value = 11690
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 7618
$ws.Cells.Item(100, 5).Value = 3900
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 172

# Row 142: Estonia - updated stats
$ws.Cells.Item(142, 2).Value = 3577
$ws.Cells.Item(142, 3).Value = 71
$ws.Cells.Item(142, 4).Value = 2727
$ws.Cells.Item(142, 5).Value = 783
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = 67
